$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 1117.1111
$ws.Range("J33").Value = 1678.6
$ws.Range("L33").Value = 1678.6
$ws.Range("N33").Value = -2136.6
# Row 42
$ws.Range("H42").Value = 187.5
$ws.Range("I42").Value = 210.75
$ws.Range("K42").Value = 632.25
$ws.Range("M42").Value = -402.25
# Row 55
$ws.Range("H55").Value = 2444
$ws.Range("I55").Value = 216.5
$ws.Range("K55").Value = 216.5
$ws.Range("M55").Value = -2.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 114
$ws.Range("H114").Value = 113000
$ws.Range("J114").Value = 113000
$ws.Range("L114").Value = 113000
$ws.Range("N114").Value = -121678
# Row 119
$ws.Range("H119").Value = 48777.668
$ws.Range("J119").Value = 48777.668
$ws.Range("L119").Value = 48777.668
$ws.Range("N119").Value = -58453.668
# Row 122
$ws.Range("H122").Value = 1638.2307
$ws.Range("I122").Value = 1461.75
$ws.Range("J122").Value = 1920.6
$ws.Range("K122").Value = 4385.25
$ws.Range("L122").Value = 5761.799999999999
$ws.Range("M122").Value = -1935.25
$ws.Range("N122").Value = -10661.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 109
$ws.Range("H109").Value = 70494.5
$ws.Range("J109").Value = 70494.5
$ws.Range("L109").Value = 70494.5
$ws.Range("N109").Value = -73268.5
# Row 127
$ws.Range("H127").Value = 60981.75
$ws.Range("J127").Value = 60981.75
$ws.Range("L127").Value = 60981.75
$ws.Range("N127").Value = -70901.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 347.5
$ws.Range("I22").Value = 347.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 347.5
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 2.5
# Row 31
$ws.Range("H31").Value = 672080.5
$ws.Range("I31").Value = 10081.294
$ws.Range("J31").Value = 1297302
$ws.Range("K31").Value = 10081.294
$ws.Range("L31").Value = 1297302
$ws.Range("M31").Value = -9786.294
$ws.Range("N31").Value = -1297892
# Row 34
$ws.Range("H34").Value = 672080.5
$ws.Range("I34").Value = 10081.294
$ws.Range("J34").Value = 1297302
$ws.Range("K34").Value = 10081.294
$ws.Range("L34").Value = 1297302
$ws.Range("M34").Value = -9879.294
$ws.Range("N34").Value = -1297706
# Row 108
$ws.Range("H108").Value = 77464.5
$ws.Range("J108").Value = 77464.5
$ws.Range("L108").Value = 77464.5
$ws.Range("N108").Value = -85144.5
# Row 132
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# Row 134
$ws.Range("H134").Value = 838421.8
$ws.Range("I134").Value = 2001415.6
$ws.Range("K134").Value = 6004246.800000001
$ws.Range("M134").Value = -6001711.800000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 297
$ws.Range("I2").Value = 414
$ws.Range("K2").Value = 2484
$ws.Range("M2").Value = -2371
# Row 10
$ws.Range("H10").Value = 25.833334
$ws.Range("I10").Value = 25.833334
$ws.Range("K10").Value = 77.50000199999999
$ws.Range("M10").Value = 61.49999800000001
# Row 23
$ws.Range("H23").Value = 872.2222
$ws.Range("J23").Value = 948.75
$ws.Range("L23").Value = 2846.25
$ws.Range("N23").Value = -3316.25
# Row 50
$ws.Range("H50").Value = 458
$ws.Range("J50").Value = 462.5
$ws.Range("L50").Value = 1387.5
$ws.Range("N50").Value = -2349.5
# Row 52
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 120000
$ws.Range("N52").Value = -120532
# Row 53
$ws.Range("H53").Value = 458
$ws.Range("J53").Value = 462.5
$ws.Range("L53").Value = 1387.5
$ws.Range("N53").Value = -2349.5
# Row 113
$ws.Range("H113").Value = 1027.0952
$ws.Range("I113").Value = 548.6667
$ws.Range("K113").Value = 1646.0001
$ws.Range("M113").Value = 523.9999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1999
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").Value = 5997
$ws.Range("N126").Value = -10937

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 5320.75
$ws.Range("I10").Value = 3765
$ws.Range("J10").Value = 9988
$ws.Range("K10").Value = 3765
$ws.Range("L10").Value = 9988
$ws.Range("M10").Value = -3625
$ws.Range("N10").Value = -10268
# Row 22
$ws.Range("H22").Value = 1662.5
$ws.Range("I22").Value = 2440.2
$ws.Range("J22").Value = 366.33334
$ws.Range("K22").Value = 2440.2
$ws.Range("L22").Value = 366.33334
$ws.Range("M22").Value = -2145.2
$ws.Range("N22").Value = -956.33334
# Row 27
$ws.Range("H27").Value = 1662.5
$ws.Range("I27").Value = 2440.2
$ws.Range("J27").Value = 366.33334
$ws.Range("K27").Value = 2440.2
$ws.Range("L27").Value = 366.33334
$ws.Range("M27").Value = -2333.2
$ws.Range("N27").Value = -580.33334
# Row 46
$ws.Range("H46").Value = 2177.8484
$ws.Range("I46").Value = 1680.3334
$ws.Range("K46").Value = 1680.3334
$ws.Range("M46").Value = -1492.3334
# Row 93
$ws.Range("H93").Value = 66668784
$ws.Range("I93").Value = 100001970
$ws.Range("K93").Value = 100001970
$ws.Range("M93").Value = -100000722
# Row 109
$ws.Range("H109").Value = 41000
$ws.Range("J109").Value = 41000
$ws.Range("L109").Value = 41000
$ws.Range("N109").Value = -43774
# Row 132
$ws.Range("H132").Value = 931957.5600000001
$ws.Range("I132").Value = 1668916.9
$ws.Range("J132").Value = 194998.33
$ws.Range("K132").Value = 5006750.699999999
$ws.Range("L132").Value = 584994.99
$ws.Range("M132").Value = -5004220.699999999
$ws.Range("N132").Value = -590054.99
# Row 136
$ws.Range("H136").Value = 354000.5
$ws.Range("I136").Value = 334667.66
$ws.Range("K136").Value = 1004002.98
$ws.Range("M136").Value = -1001452.98

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 27247.5
$ws.Range("I40").Value = 24000
$ws.Range("J40").Value = 30495
$ws.Range("K40").Value = 24000
$ws.Range("L40").Value = 30495
$ws.Range("M40").Value = -23851
$ws.Range("N40").Value = -30793
# Row 81
$ws.Range("H81").Value = 1236.5
$ws.Range("J81").Value = 1950
$ws.Range("L81").Value = 3900
$ws.Range("N81").Value = -6022
# Row 84
$ws.Range("H84").Value = 1236.5
$ws.Range("J84").Value = 1950
$ws.Range("L84").Value = 19500
$ws.Range("N84").Value = -30108
# Row 107
$ws.Range("H107").Value = 31251696
$ws.Range("J107").Value = 922.6
$ws.Range("L107").Value = 2767.8
$ws.Range("N107").Value = -6607.8
# Row 119
$ws.Range("H119").Value = 61660
$ws.Range("J119").Value = 61660
$ws.Range("L119").Value = 61660
$ws.Range("N119").Value = -71336
# Row 123
$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -69800
# Row 132
$ws.Range("H132").Value = 24064.777
$ws.Range("I132").Value = 2527.6667
$ws.Range("K132").Value = 7583.000100000001
$ws.Range("M132").Value = -5053.000100000001
# Row 136
$ws.Range("H136").Value = 34000
$ws.Range("J136").Value = 38571.43
$ws.Range("L136").Value = 115714.29
$ws.Range("N136").Value = -120814.29
